$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Title rows (D1:D4) -------------------------------------------------
$ws.Range("D1").Value = "NR Finance Mexico"
$ws.Range("D2").Value = "EKIP"
$ws.Range("D3").Value = "Certificacion de usuarios 2024"
$ws.Range("D4").Value = "Reporte de usuarios"

$titles = $ws.Range("D1:D4")
$titles.Font.Name = "Calibri"
$titles.Font.Size = 16
$titles.Font.Bold = $true
$titles.Font.Color = 0
$titles.HorizontalAlignment = -4108   # xlCenter

# --- Data table (A5:F11) : thin border box around every cell -----------
$table = $ws.Range("A5:F11")
$table.Borders.LineStyle = 1          # xlContinuous
$table.Borders.Weight = 2             # xlThin

# --- Column widths -------------------------------------------------------
# Excel's ColumnWidth setter stores the width rounded to the nearest 1/6th
# of a character (pixel-snapped); feed it (target - 5/6) so the persisted
# width lands on the closest achievable value to the desired target.
$ws.Columns.Item(1).ColumnWidth = 4.853482 - (5/6)
$ws.Columns.Item(2).ColumnWidth = 30.567768 - (5/6)
$ws.Columns.Item(3).ColumnWidth = 8.424911 - (5/6)
$ws.Columns.Item(4).ColumnWidth = 40.424911 - (5/6)
$ws.Columns.Item(5).ColumnWidth = 24.139196 - (5/6)
$ws.Columns.Item(6).ColumnWidth = 14.996339 - (5/6)
